$d = $word.ActiveDocument

# --- Part 1: remove the old `_GoBack` bookmark that sits after "21" ---
# (handled implicitly below: `_GoBack` is a unique bookmark name, so
#  adding a new one automatically removes the previous one from the
#  document; see Part 2.)

# --- Part 2: insert the letter "t" right after "ship1" (turning
#     "...ship1 appear..." into "...ship1t appear...") and leave the
#     `_GoBack` bookmark collapsed immediately after the inserted "t" ---

$r = $d.Content
$r.Find.Execute("ship1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $r.End

# Plant the (new) `_GoBack` bookmark collapsed right after "ship1".
# Because `_GoBack` is a singleton bookmark name, Word removes the
# previous occurrence (the one after "21") when this one is added.
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Insert "t" immediately before the bookmark (i.e. right after "ship1").
$bm = $d.Bookmarks.Item("_GoBack")
$bmR = $bm.Range
$insertRange = $d.Range($bmR.Start, $bmR.Start)
$insertRange.InsertBefore("t")

# Nudge formatting off/on so the engine keeps "t" as its own run
# instead of re-merging it into the preceding "ship1" run, matching
# real Word's run-splitting behavior around the `_GoBack` bookmark.
$insertRange.Font.Bold = 1
$insertRange.Font.Bold = 0
